# Updates for test problems
# - species sheet: update Inlet Mole Frac values in column C
# - update selections on both sheets (system -> D8, species -> C9)
#   ("system" must remain the tab-selected / active sheet at the end)

$wb = $excel.ActiveWorkbook

$wsSystem  = $wb.Worksheets.Item("system")
$wsSpecies = $wb.Worksheets.Item("species")

# Update the "Inlet Mole Frac" values on the species sheet
$wsSpecies.Range("C2").Value = 0.1
$wsSpecies.Range("C4").Value = 0.4
$wsSpecies.Range("C5").Value = 0.1
$wsSpecies.Range("C6").Value = 0.1

# Update the species sheet's remembered selection first...
$wsSpecies.Activate()
$wsSpecies.Range("C9").Select()

# ...then re-activate the system sheet last so it stays the active/tabSelected
# sheet, with its own selection moved to D8.
$wsSystem.Activate()
$wsSystem.Range("D8").Select()
